$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting of the last existing data row (104) down onto the two new rows
$ws.Range("A104:E104").Copy()
$ws.Range("A105:E105").PasteSpecial(-4122)
$ws.Range("A104:E104").Copy()
$ws.Range("A106:E106").PasteSpecial(-4122)

# Fill in the new test case data (order chosen to match shared-string build order)
$ws.Range("A105").Value = "TestCase_B104"
$ws.Range("A106").Value = "TestCase_B105"

$ws.Range("B105").Value = "OPQA-554"
$ws.Range("B106").Value = "OPQA-555"

$ws.Range("C106").Value = "Verify that record view page of a post gets displayed when user clicks on article title in POSTs search results page"
$ws.Range("C105").Value = "Verify that record view page of a post gets displayed when user clicks on article title in ALL  search results page"

$ws.Range("D105").Value = "Y"
$ws.Range("D106").Value = "Y"

# Restore the selection to where the author ended up
$ws.Range("C97").Select() | Out-Null
